$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H64").Value = 65545.06
$ws_ALC.Range("J64").Value = 3412.1
$ws_ALC.Range("L64").Value = 3412.1
$ws_ALC.Range("N64").Value = -3908.1
$ws_ALC.Range("H67").Value = 65545.06
$ws_ALC.Range("J67").Value = 3412.1
$ws_ALC.Range("L67").Value = 3412.1
$ws_ALC.Range("N67").Value = -5128.1
$ws_ALC.Range("H112").Value = 1563.9584
$ws_ALC.Range("J112").Value = 1563.9584
$ws_ALC.Range("L112").Value = 4691.8752
$ws_ALC.Range("N112").Value = -6907.8752
$ws_ALC.Range("H118").Value = 6331.6113
$ws_ALC.Range("I118").Value = 9427.272000000001
$ws_ALC.Range("J118").Value = 1467
$ws_ALC.Range("K118").Value = 28281.816
$ws_ALC.Range("L118").Value = 4401
$ws_ALC.Range("M118").Value = -26624.816
$ws_ALC.Range("N118").Value = -7715
$ws_ALC.Range("H132").Value = 3792823
$ws_ALC.Range("I132").Value = 4171538.8
$ws_ALC.Range("K132").Value = 12514616.4
$ws_ALC.Range("M132").Value = -12512086.4

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 6267.64
$ws_ARM.Range("I32").Value = 6289.156
$ws_ARM.Range("J32").Value = 6195.609
$ws_ARM.Range("K32").Value = 6289.156
$ws_ARM.Range("L32").Value = 6195.609
$ws_ARM.Range("M32").Value = -6002.156
$ws_ARM.Range("N32").Value = -6769.609
$ws_ARM.Range("H74").Value = 754.8611
$ws_ARM.Range("I74").Value = 623.5185
$ws_ARM.Range("J74").Value = 1148.8889
$ws_ARM.Range("K74").Value = 623.5185
$ws_ARM.Range("L74").Value = 1148.8889
$ws_ARM.Range("M74").Value = 250.4815
$ws_ARM.Range("N74").Value = -2896.8889
$ws_ARM.Range("H77").Value = 754.8611
$ws_ARM.Range("I77").Value = 623.5185
$ws_ARM.Range("J77").Value = 1148.8889
$ws_ARM.Range("K77").Value = 3117.5925
$ws_ARM.Range("L77").Value = 5744.4445
$ws_ARM.Range("M77").Value = 1250.4075
$ws_ARM.Range("N77").Value = -14480.4445

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H26").Value = 10111.833
$ws_BSM.Range("I26").Value = 10111.833
$ws_BSM.Range("J26").Value = 0
$ws_BSM.Range("K26").Value = 10111.833
$ws_BSM.Range("L26").Value = 0
$ws_BSM.Range("M26").Value = -9819.833000000001
$ws_BSM.Range("N26").ClearContents()
$ws_BSM.Range("H96").Value = 10332.8
$ws_BSM.Range("I96").Value = 5691
$ws_BSM.Range("J96").Value = 28900
$ws_BSM.Range("K96").Value = 5691
$ws_BSM.Range("L96").Value = 28900
$ws_BSM.Range("M96").Value = -2945
$ws_BSM.Range("N96").Value = -34392
$ws_BSM.Range("H105").Value = 112912.61
$ws_BSM.Range("I105").Value = 85153.836
$ws_BSM.Range("J105").Value = 168430.17
$ws_BSM.Range("K105").Value = 85153.836
$ws_BSM.Range("L105").Value = 168430.17
$ws_BSM.Range("M105").Value = -83406.836
$ws_BSM.Range("N105").Value = -171924.17

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 29266.309
$ws_CRP.Range("I31").Value = 889.5
$ws_CRP.Range("K31").Value = 889.5
$ws_CRP.Range("M31").Value = -594.5
$ws_CRP.Range("H34").Value = 29266.309
$ws_CRP.Range("I34").Value = 889.5
$ws_CRP.Range("K34").Value = 889.5
$ws_CRP.Range("M34").Value = -687.5
$ws_CRP.Range("H99").Value = 2220.6667
$ws_CRP.Range("I99").Value = 1960.8889
$ws_CRP.Range("K99").Value = 1960.8889
$ws_CRP.Range("M99").Value = -462.8888999999999
$ws_CRP.Range("H126").Value = 2220.6667
$ws_CRP.Range("I126").Value = 1960.8889
$ws_CRP.Range("K126").Value = 5882.6667
$ws_CRP.Range("M126").Value = -3412.6667
$ws_CRP.Range("H132").Value = 6400.222
$ws_CRP.Range("I132").Value = 7634
$ws_CRP.Range("J132").Value = 3932.6667
$ws_CRP.Range("K132").Value = 22902
$ws_CRP.Range("L132").Value = 11798.0001
$ws_CRP.Range("M132").Value = -20372
$ws_CRP.Range("N132").Value = -16858.0001
$ws_CRP.Range("H134").Value = 1032.2222
$ws_CRP.Range("I134").Value = 557.62067
$ws_CRP.Range("J134").Value = 2998.4285
$ws_CRP.Range("K134").Value = 1672.86201
$ws_CRP.Range("L134").Value = 8995.2855
$ws_CRP.Range("M134").Value = 862.1379899999999
$ws_CRP.Range("N134").Value = -14065.2855

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H34").Value = 848.7778
$ws_CUL.Range("I34").Value = 25
$ws_CUL.Range("J34").Value = 1084.1428
$ws_CUL.Range("K34").Value = 75
$ws_CUL.Range("L34").Value = 3252.4284
$ws_CUL.Range("M34").Value = 9
$ws_CUL.Range("N34").Value = -3420.4284
$ws_CUL.Range("H107").Value = 969.1667
$ws_CUL.Range("I107").Value = 561.3333
$ws_CUL.Range("J107").Value = 1260.4762
$ws_CUL.Range("K107").Value = 1683.9999
$ws_CUL.Range("L107").Value = 3781.4286
$ws_CUL.Range("M107").Value = 236.0001
$ws_CUL.Range("N107").Value = -7621.4286
$ws_CUL.Range("H113").Value = 944.25
$ws_CUL.Range("I113").Value = 1273.6154
$ws_CUL.Range("J113").Value = 658.8
$ws_CUL.Range("K113").Value = 3820.8462
$ws_CUL.Range("L113").Value = 1976.4
$ws_CUL.Range("M113").Value = -1650.8462
$ws_CUL.Range("N113").Value = -6316.4
$ws_CUL.Range("H131").Value = 507331.1
$ws_CUL.Range("I131").Value = 777.9375
$ws_CUL.Range("J131").Value = 603817.4
$ws_CUL.Range("K131").Value = 2333.8125
$ws_CUL.Range("L131").Value = 1811452.2
$ws_CUL.Range("M131").Value = 2706.1875
$ws_CUL.Range("N131").Value = -1821532.2
$ws_CUL.Range("H140").Value = 5404.8
$ws_CUL.Range("I140").Value = 5835
$ws_CUL.Range("K140").Value = 17505
$ws_CUL.Range("M140").Value = -12325

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H43").Value = 2033.4828
$ws_GSM.Range("I43").Value = 998
$ws_GSM.Range("J43").Value = 7003.8
$ws_GSM.Range("K43").Value = 998
$ws_GSM.Range("L43").Value = 7003.8
$ws_GSM.Range("M43").Value = -847
$ws_GSM.Range("N43").Value = -7305.8
$ws_GSM.Range("H46").Value = 11966.533
$ws_GSM.Range("I46").Value = 5000
$ws_GSM.Range("J46").Value = 12464.143
$ws_GSM.Range("K46").Value = 5000
$ws_GSM.Range("L46").Value = 12464.143
$ws_GSM.Range("M46").Value = -4844
$ws_GSM.Range("N46").Value = -12776.143
$ws_GSM.Range("H57").Value = 5000
$ws_GSM.Range("I57").Value = 5000
$ws_GSM.Range("K57").Value = 5000
$ws_GSM.Range("M57").Value = -4180
$ws_GSM.Range("H126").Value = 3679482
$ws_GSM.Range("I126").Value = 3200.3635
$ws_GSM.Range("J126").Value = 11767301
$ws_GSM.Range("K126").Value = 9601.0905
$ws_GSM.Range("L126").Value = 35301903
$ws_GSM.Range("M126").Value = -7131.0905
$ws_GSM.Range("N126").Value = -35306843

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 2515.6
$ws_LTW.Range("I7").Value = 1659.3334
$ws_LTW.Range("J7").Value = 3800
$ws_LTW.Range("K7").Value = 1659.3334
$ws_LTW.Range("L7").Value = 3800
$ws_LTW.Range("M7").Value = -1547.3334
$ws_LTW.Range("N7").Value = -4024
$ws_LTW.Range("H40").Value = 58640.89
$ws_LTW.Range("I40").Value = 146155.14
$ws_LTW.Range("K40").Value = 146155.14
$ws_LTW.Range("M40").Value = -146019.14
$ws_LTW.Range("H68").Value = 3926.2856
$ws_LTW.Range("I68").Value = 2216.8333
$ws_LTW.Range("K68").Value = 2216.8333
$ws_LTW.Range("M68").Value = -1467.8333
$ws_LTW.Range("H71").Value = 3926.2856
$ws_LTW.Range("I71").Value = 2216.8333
$ws_LTW.Range("K71").Value = 11084.1665
$ws_LTW.Range("M71").Value = -7340.166499999999
$ws_LTW.Range("H126").Value = 2515.6
$ws_LTW.Range("I126").Value = 1659.3334
$ws_LTW.Range("J126").Value = 3800
$ws_LTW.Range("K126").Value = 4978.0002
$ws_LTW.Range("L126").Value = 11400
$ws_LTW.Range("M126").Value = -2508.0002
$ws_LTW.Range("N126").Value = -16340
$ws_LTW.Range("H136").Value = 1963.0975
$ws_LTW.Range("I136").Value = 1387.9333
$ws_LTW.Range("J136").Value = 3531.7273
$ws_LTW.Range("K136").Value = 4163.7999
$ws_LTW.Range("L136").Value = 10595.1819
$ws_LTW.Range("M136").Value = -1613.7999
$ws_LTW.Range("N136").Value = -15695.1819

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H116").Value = 43000
$ws_WVR.Range("J116").Value = 43000
$ws_WVR.Range("L116").Value = 43000
$ws_WVR.Range("N116").Value = -52178
$ws_WVR.Range("H124").Value = 28714.223
$ws_WVR.Range("J124").Value = 28714.223
$ws_WVR.Range("L124").Value = 28714.223
$ws_WVR.Range("N124").Value = -38534.223
